{"js": "// Update the \"Do B\u00e0 / NGUY\u1ec4N TH\u1eca B\u00cdCH NG\u00c2N / Ph\u00f3 Gi\u00e1m \u0111\u1ed1c\" signatory block\n// to the new \"Do \u00d4ng / NGUY\u1ec4N QU\u1ed0C \u0110\u1ea0T / Tr\u01b0\u1edfng ph\u00f2ng kinh doanh\" signatory\n// (template update for Mr. Dat), and move the `_GoBack` bookmark from the\n// \"T\u00e0i kho\u1ea3n (4)\" paragraph to right after the new job title.\n\nconst body = context.document.body;\n\n// Locate the signatory paragraph (\"Do B\u00e0 ... NGUY\u1ec4N TH\u1eca B\u00cdCH NG\u00c2N ... Ch\u1ee9c v\u1ee5: ...\").\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet signatoryIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"NGUY\u1ec4N TH\u1eca B\u00cdCH NG\u00c2N\") !== -1) {\n    signatoryIndex = i;\n    break;\n  }\n}\nif (signatoryIndex === -1) {\n  throw new Error(\"Could not find the signatory paragraph\");\n}\nconst signatoryParagraph = paragraphs.items[signatoryIndex];\n\n// 1) \"Do B\u00e0 \" -> \"Do \u00d4ng \"\n{\n  const found = signatoryParagraph.search(\"B\u00e0 \", { matchCase: true, matchWholeWord: false });\n  found.load(\"text\");\n  await context.sync();\n  if (found.items.length > 0) {\n    found.items[0].insertText(\"\u00d4ng \", Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n\n// 2) \"NGUY\u1ec4N TH\u1eca B\u00cdCH NG\u00c2N\" -> \"NGUY\u1ec4N QU\u1ed0C \u0110\u1ea0T\" (keeps bold formatting of the run)\n{\n  const found = signatoryParagraph.search(\"NGUY\u1ec4N TH\u1eca B\u00cdCH NG\u00c2N\", { matchCase: true });\n  found.load(\"text\");\n  await context.sync();\n  if (found.items.length > 0) {\n    found.items[0].insertText(\"NGUY\u1ec4N QU\u1ed0C \u0110\u1ea0T\", Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n\n// 3) \"Ch\u1ee9c v\u1ee5:  Ph\u00f3 Gi\u00e1m \u0111\u1ed1c  \u2013  l\u00e0m \u0111\u1ea1i di\u1ec7n.\" -> \"Ch\u1ee9c v\u1ee5:  Tr\u01b0\u1edfng ph\u00f2ng kinh doanh.\"\n{\n  const found = signatoryParagraph.search(\"Ph\u00f3 Gi\u00e1m \u0111\u1ed1c  \u2013  l\u00e0m \u0111\u1ea1i di\u1ec7n.\", { matchCase: true });\n  found.load(\"text\");\n  await context.sync();\n  if (found.items.length > 0) {\n    found.items[0].insertText(\"Tr\u01b0\u1edfng ph\u00f2ng kinh doanh.\", Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n\n// 4) Move the `_GoBack` bookmark from its old spot (\"T\u00e0i kho\u1ea3n (4)\" paragraph)\n//    to just after the new job title \"Tr\u01b0\u1edfng ph\u00f2ng kinh doanh\" (before the final \".\").\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n{\n  const refreshedParagraphs = body.paragraphs;\n  refreshedParagraphs.load(\"text\");\n  await context.sync();\n\n  let idx = -1;\n  for (let i = 0; i < refreshedParagraphs.items.length; i++) {\n    if (refreshedParagraphs.items[i].text.indexOf(\"Tr\u01b0\u1edfng ph\u00f2ng kinh doanh\") !== -1) {\n      idx = i;\n      break;\n    }\n  }\n  if (idx === -1) {\n    throw new Error(\"Could not find the updated job-title paragraph\");\n  }\n  const titleParagraph = refreshedParagraphs.items[idx];\n\n  const found = titleParagraph.search(\"Tr\u01b0\u1edfng ph\u00f2ng kinh doanh\", { matchCase: true });\n  found.load(\"text\");\n  await context.sync();\n  if (found.items.length > 0) {\n    const endRange = found.items[0].getRange(\"End\");\n    endRange.insertBookmark(\"_GoBack\");\n    await context.sync();\n  }\n}\n", "ps1": "# Update the \"Do B\u00e0 / NGUY\u1ec4N TH\u1eca B\u00cdCH NG\u00c2N / Ph\u00f3 Gi\u00e1m \u0111\u1ed1c\" signatory block\n# to the new \"Do \u00d4ng / NGUY\u1ec4N QU\u1ed0C \u0110\u1ea0T / Tr\u01b0\u1edfng ph\u00f2ng kinh doanh\" signatory\n# (template update for Mr. Dat), and move the `_GoBack` bookmark from the\n# \"T\u00e0i kho\u1ea3n (4)\" paragraph to right after the new job title.\n\n$d = $word.ActiveDocument\n\nfunction Find-ParagraphContaining([string]$needle) {\n    $r = $d.Content\n    $f = $r.Find\n    $f.Text = $needle\n    $f.MatchCase = $true\n    $ok = $f.Execute()\n    if (-not $ok) {\n        throw \"Could not find paragraph containing: $needle\"\n    }\n    return $r.Paragraphs(1).Range\n}\n\n# 1) \"Do B\u00e0 \" -> \"Do \u00d4ng \"\n$p1 = Find-ParagraphContaining \"NGUY\u1ec4N TH\u1eca B\u00cdCH NG\u00c2N\"\n$p1.Find.Execute(\"B\u00e0 \", $true, $false, $false, $false, $false, $true, 0, $false, \"\u00d4ng \", 1) | Out-Null\n\n# 2) \"NGUY\u1ec4N TH\u1eca B\u00cdCH NG\u00c2N\" -> \"NGUY\u1ec4N QU\u1ed0C \u0110\u1ea0T\" (keeps the bold run formatting)\n$p2 = Find-ParagraphContaining \"NGUY\u1ec4N TH\u1eca B\u00cdCH NG\u00c2N\"\n$p2.Find.Execute(\"NGUY\u1ec4N TH\u1eca B\u00cdCH NG\u00c2N\", $true, $false, $false, $false, $false, $true, 0, $false, \"NGUY\u1ec4N QU\u1ed0C \u0110\u1ea0T\", 1) | Out-Null\n\n# 3) \"Ch\u1ee9c v\u1ee5:  Ph\u00f3 Gi\u00e1m \u0111\u1ed1c  \u2013  l\u00e0m \u0111\u1ea1i di\u1ec7n.\" -> \"Ch\u1ee9c v\u1ee5:  Tr\u01b0\u1edfng ph\u00f2ng kinh doanh.\"\n$p3 = Find-ParagraphContaining \"NGUY\u1ec4N QU\u1ed0C \u0110\u1ea0T\"\n$p3.Find.Execute(\"Ph\u00f3 Gi\u00e1m \u0111\u1ed1c  \u2013  l\u00e0m \u0111\u1ea1i di\u1ec7n.\", $true, $false, $false, $false, $false, $true, 0, $false, \"Tr\u01b0\u1edfng ph\u00f2ng kinh doanh.\", 1) | Out-Null\n\n# 4) Move the `_GoBack` bookmark from its old spot (\"T\u00e0i kho\u1ea3n (4)\" paragraph)\n#    to just after the new job title \"Tr\u01b0\u1edfng ph\u00f2ng kinh doanh\" (before the final \".\").\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n$p4 = Find-ParagraphContaining \"Tr\u01b0\u1edfng ph\u00f2ng kinh doanh\"\n$p4.Find.Execute(\"Tr\u01b0\u1edfng ph\u00f2ng kinh doanh\", $true, $false, $false, $false, $false, $true, 0) | Out-Null\n$bookmarkRange = $p4.Duplicate\n$bookmarkRange.Collapse(0)\n$d.Bookmarks.Add(\"_GoBack\", $bookmarkRange) | Out-Null\n"}
